$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 399
$ws.Range("I2").Value = 1042
$ws.Range("J2").Value = 4139
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 1137
$ws.Range("M2").Value = 58
$ws.Range("N2").Value = 698
$ws.Range("P2").Value = 19
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 43
$ws.Range("S2").Value = 414
$ws.Range("T2").Value = 711
$ws.Range("U2").Value = 52
$ws.Range("V2").Value = 6351
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 6478
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 79
$ws.Range("AA2").Value = 34
